$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change from numeric values to inline strings
$ws.Range("A2").Value = "Pd"
$ws.Range("B2").Value = "Cd"
$ws.Range("C2").Value = "Ru"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 21

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "#N/A"
$ws.Range("C4").Value = 22

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = 23
